$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set cell values
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Style B1: bold font, thin box border, centered horizontally, top vertically
$r1 = $ws.Range("B1")
$r1.Borders.LineStyle = 1    # xlContinuous
$r1.Borders.Weight = 2       # xlThin
$r1.Font.Bold = $true
$r1.HorizontalAlignment = -4108  # xlCenter
$r1.VerticalAlignment = -4160    # xlTop

# Copy B1's format (and value) onto A2 so both share the same style entry
# instead of rebuilding it property-by-property (which would leave unused
# intermediate style table entries behind).
$r1.Copy($ws.Range("A2"))
